$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header row with a new "time_taken" column (F1), copying the
# formatting (bold, centered, bordered) from the existing header cell E1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Per-row timestamps recorded for column F (rows 2-100), taken verbatim
# from the panelapp scrape metadata.
$timeTaken = @(
    "2021-10-05 13:38:52.234855",
    "2021-10-05 13:38:52.234864",
    "2021-10-05 13:38:52.234867",
    "2021-10-05 13:38:52.234869",
    "2021-10-05 13:38:52.234871",
    "2021-10-05 13:38:52.234873",
    "2021-10-05 13:38:52.234875",
    "2021-10-05 13:38:52.234877",
    "2021-10-05 13:38:52.234879",
    "2021-10-05 13:38:52.234881",
    "2021-10-05 13:38:52.234883",
    "2021-10-05 13:38:52.234885",
    "2021-10-05 13:38:52.234887",
    "2021-10-05 13:38:52.234889",
    "2021-10-05 13:38:52.234891",
    "2021-10-05 13:38:52.234893",
    "2021-10-05 13:38:52.234895",
    "2021-10-05 13:38:52.234897",
    "2021-10-05 13:38:52.234899",
    "2021-10-05 13:38:52.234901",
    "2021-10-05 13:38:52.234903",
    "2021-10-05 13:38:52.234905",
    "2021-10-05 13:38:52.234907",
    "2021-10-05 13:38:52.234909",
    "2021-10-05 13:38:52.234911",
    "2021-10-05 13:38:52.234913",
    "2021-10-05 13:38:52.234915",
    "2021-10-05 13:38:52.234917",
    "2021-10-05 13:38:52.234918",
    "2021-10-05 13:38:52.234920",
    "2021-10-05 13:38:52.234922",
    "2021-10-05 13:38:52.234924",
    "2021-10-05 13:38:52.234926",
    "2021-10-05 13:38:52.234928",
    "2021-10-05 13:38:52.234930",
    "2021-10-05 13:38:52.234932",
    "2021-10-05 13:38:52.234934",
    "2021-10-05 13:38:52.234936",
    "2021-10-05 13:38:52.234938",
    "2021-10-05 13:38:52.234940",
    "2021-10-05 13:38:52.234942",
    "2021-10-05 13:38:52.234944",
    "2021-10-05 13:38:52.234946",
    "2021-10-05 13:38:52.234948",
    "2021-10-05 13:38:52.234950",
    "2021-10-05 13:38:52.234952",
    "2021-10-05 13:38:52.234954",
    "2021-10-05 13:38:52.234955",
    "2021-10-05 13:38:52.234957",
    "2021-10-05 13:38:52.234959",
    "2021-10-05 13:38:52.234961",
    "2021-10-05 13:38:52.234963",
    "2021-10-05 13:38:52.234965",
    "2021-10-05 13:38:52.234967",
    "2021-10-05 13:38:52.234969",
    "2021-10-05 13:38:52.234971",
    "2021-10-05 13:38:52.234973",
    "2021-10-05 13:38:52.234975",
    "2021-10-05 13:38:52.234977",
    "2021-10-05 13:38:52.234979",
    "2021-10-05 13:38:52.234981",
    "2021-10-05 13:38:52.234983",
    "2021-10-05 13:38:52.234985",
    "2021-10-05 13:38:52.234987",
    "2021-10-05 13:38:52.234990",
    "2021-10-05 13:38:52.234992",
    "2021-10-05 13:38:52.234994",
    "2021-10-05 13:38:52.234996",
    "2021-10-05 13:38:52.234998",
    "2021-10-05 13:38:52.235000",
    "2021-10-05 13:38:52.235002",
    "2021-10-05 13:38:52.235004",
    "2021-10-05 13:38:52.235006",
    "2021-10-05 13:38:52.235008",
    "2021-10-05 13:38:52.235010",
    "2021-10-05 13:38:52.235012",
    "2021-10-05 13:38:52.235015",
    "2021-10-05 13:38:52.235018",
    "2021-10-05 13:38:52.235019",
    "2021-10-05 13:38:52.235021",
    "2021-10-05 13:38:52.235023",
    "2021-10-05 13:38:52.235025",
    "2021-10-05 13:38:52.235027",
    "2021-10-05 13:38:52.235029",
    "2021-10-05 13:38:52.235032",
    "2021-10-05 13:38:52.235033",
    "2021-10-05 13:38:52.235035",
    "2021-10-05 13:38:52.235037",
    "2021-10-05 13:38:52.235039",
    "2021-10-05 13:38:52.235041",
    "2021-10-05 13:38:52.235043",
    "2021-10-05 13:38:52.235045",
    "2021-10-05 13:38:52.235048",
    "2021-10-05 13:38:52.235050",
    "2021-10-05 13:38:52.235052",
    "2021-10-05 13:38:52.235054",
    "2021-10-05 13:38:52.235056",
    "2021-10-05 13:38:52.235058",
    "2021-10-05 13:38:52.235060"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

Write-Output ("Wrote time_taken header plus {0} timestamp rows" -f $timeTaken.Length)
